$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old sample data row (row 2)
$ws.Rows("2:2").Delete()

# Insert a new column for "Modes de paiement" right after "Cause de l'incident" (column C),
# shifting Site..Initier par (D:Q) one column to the right (E:R)
$ws.Columns("D:D").Insert()
$ws.Cells.Item(1, 4).Value = "Modes de paiement "
